$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2080.5833
$ws.Range("I15").Value = 2080.5833
$ws.Range("K15").Value = 6241.749899999999
$ws.Range("M15").Value = -6072.749899999999

$ws.Range("H40").Value = 2327.6155
$ws.Range("J40").Value = 2298.8572
$ws.Range("L40").Value = 2298.8572
$ws.Range("N40").Value = -2648.8572

$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 7000
$ws.Range("K62").Value = 7000
$ws.Range("M62").Value = -6376

$ws.Range("H64").Value = 5365.5
$ws.Range("I64").Value = 4299.25
$ws.Range("K64").Value = 4299.25
$ws.Range("M64").Value = -4051.25

$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 7000
$ws.Range("K65").Value = 35000
$ws.Range("M65").Value = -31880

$ws.Range("H67").Value = 5365.5
$ws.Range("I67").Value = 4299.25
$ws.Range("K67").Value = 4299.25
$ws.Range("M67").Value = -3441.25

$ws.Range("H80").Value = 2178.35
$ws.Range("I80").Value = 2133.5
$ws.Range("J80").Value = 2283
$ws.Range("K80").Value = 6400.5
$ws.Range("L80").Value = 6849
$ws.Range("M80").Value = -5402.5
$ws.Range("N80").Value = -8845

$ws.Range("H83").Value = 2178.35
$ws.Range("I83").Value = 2133.5
$ws.Range("J83").Value = 2283
$ws.Range("K83").Value = 19201.5
$ws.Range("L83").Value = 20547
$ws.Range("M83").Value = -14209.5
$ws.Range("N83").Value = -30531

$ws.Range("H113").Value = 5431.5
$ws.Range("I113").Value = 3343.8
$ws.Range("J113").Value = 7519.2
$ws.Range("K113").Value = 3343.8
$ws.Range("L113").Value = 7519.2
$ws.Range("M113").Value = -89.80000000000018
$ws.Range("N113").Value = -14027.2

$ws.Range("H129").Value = 2346.2727
$ws.Range("I129").Value = 2175
$ws.Range("J129").Value = 2363.4
$ws.Range("K129").Value = 6525
$ws.Range("L129").Value = 7090.200000000001
$ws.Range("M129").Value = -1525
$ws.Range("N129").Value = -17090.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1504.5
$ws.Range("I2").Value = 1416.6666
$ws.Range("J2").Value = 1899.75
$ws.Range("K2").Value = 1416.6666
$ws.Range("L2").Value = 1899.75
$ws.Range("M2").Value = -1303.6666
$ws.Range("N2").Value = -2125.75

$ws.Range("H45").Value = 1003.3333
$ws.Range("I45").Value = 1003.3333
$ws.Range("K45").Value = 1003.3333
$ws.Range("M45").Value = -626.3333

$ws.Range("H61").Value = 2903.6428
$ws.Range("I61").Value = 2903.6428
$ws.Range("K61").Value = 2903.6428
$ws.Range("M61").Value = -2691.6428

$ws.Range("H63").Value = 3268.756
$ws.Range("I63").Value = 1972.1143
$ws.Range("J63").Value = 10832.5
$ws.Range("K63").Value = 1972.1143
$ws.Range("L63").Value = 10832.5
$ws.Range("M63").Value = -1286.1143
$ws.Range("N63").Value = -12204.5

$ws.Range("H66").Value = 3268.756
$ws.Range("I66").Value = 1972.1143
$ws.Range("J66").Value = 10832.5
$ws.Range("K66").Value = 9860.5715
$ws.Range("L66").Value = 54162.5
$ws.Range("M66").Value = -6428.5715
$ws.Range("N66").Value = -61026.5

$ws.Range("H74").Value = 3331
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 3331
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H97").Value = 501.5
$ws.Range("I97").Value = 501.5
$ws.Range("K97").Value = 501.5
$ws.Range("M97").Value = -5.5

$ws.Range("H102").Value = 1521.0769
$ws.Range("I102").Value = 1521.0769
$ws.Range("K102").Value = 1521.0769
$ws.Range("M102").Value = 100.9231

$ws.Range("H116").Value = 1504.5
$ws.Range("I116").Value = 1416.6666
$ws.Range("J116").Value = 1899.75
$ws.Range("K116").Value = 1416.6666
$ws.Range("L116").Value = 1899.75
$ws.Range("M116").Value = 877.3334
$ws.Range("N116").Value = -6487.75

$ws.Range("H132").Value = 1550.2188
$ws.Range("I132").Value = 1471.5
$ws.Range("J132").Value = 1786.375
$ws.Range("K132").Value = 4414.5
$ws.Range("L132").Value = 5359.125
$ws.Range("M132").Value = -1884.5
$ws.Range("N132").Value = -10419.125

$ws.Range("H136").Value = 2903.6428
$ws.Range("I136").Value = 2903.6428
$ws.Range("K136").Value = 8710.928400000001
$ws.Range("M136").Value = -6160.928400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1504.5
$ws.Range("I3").Value = 1416.6666
$ws.Range("J3").Value = 1899.75
$ws.Range("K3").Value = 1416.6666
$ws.Range("L3").Value = 1899.75
$ws.Range("M3").Value = -1302.6666
$ws.Range("N3").Value = -2127.75

$ws.Range("H20").Value = 1062.4
$ws.Range("I20").Value = 970.6667
$ws.Range("J20").Value = 1200
$ws.Range("K20").Value = 970.6667
$ws.Range("L20").Value = 1200
$ws.Range("M20").Value = -723.6667
$ws.Range("N20").Value = -1694

$ws.Range("H22").Value = 258.63635
$ws.Range("I22").Value = 271.33334
$ws.Range("K22").Value = 271.33334
$ws.Range("M22").Value = -98.33334000000002

$ws.Range("H134").Value = 3441.375
$ws.Range("I134").Value = 3441.375
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10324.125
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7789.125
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4191.636
$ws.Range("I31").Value = 3345.4443
$ws.Range("J31").Value = 7999.5
$ws.Range("K31").Value = 3345.4443
$ws.Range("L31").Value = 7999.5
$ws.Range("M31").Value = -3050.4443
$ws.Range("N31").Value = -8589.5

$ws.Range("H34").Value = 4191.636
$ws.Range("I34").Value = 3345.4443
$ws.Range("J34").Value = 7999.5
$ws.Range("K34").Value = 3345.4443
$ws.Range("L34").Value = 7999.5
$ws.Range("M34").Value = -3143.4443
$ws.Range("N34").Value = -8403.5

$ws.Range("H59").Value = 66000
$ws.Range("I59").Value = 52500
$ws.Range("J59").Value = 75000
$ws.Range("K59").Value = 52500
$ws.Range("L59").Value = 75000
$ws.Range("M59").Value = -51355
$ws.Range("N59").Value = -77290

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H107").Value = 1403.8
$ws.Range("I107").Value = 1377.5
$ws.Range("J107").Value = 1465.1666
$ws.Range("K107").Value = 1377.5
$ws.Range("L107").Value = 1465.1666
$ws.Range("M107").Value = 542.5
$ws.Range("N107").Value = -5305.1666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 7925.316
$ws.Range("I9").Value = 895
$ws.Range("J9").Value = 9800.066000000001
$ws.Range("K9").Value = 2685
$ws.Range("L9").Value = 29400.198
$ws.Range("M9").Value = -2461
$ws.Range("N9").Value = -29848.198

$ws.Range("H36").Value = 2983
$ws.Range("I36").Value = 2474
$ws.Range("J36").Value = 3237.5
$ws.Range("K36").Value = 7422
$ws.Range("L36").Value = 9712.5
$ws.Range("M36").Value = -7253
$ws.Range("N36").Value = -10050.5

$ws.Range("H131").Value = 2178.8572
$ws.Range("I131").Value = 1963
$ws.Range("J131").Value = 2466.6667
$ws.Range("K131").Value = 5889
$ws.Range("L131").Value = 7400.000100000001
$ws.Range("M131").Value = -849
$ws.Range("N131").Value = -17480.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 11982.333
$ws.Range("I46").Value = 2950.5
$ws.Range("J46").Value = 30046
$ws.Range("K46").Value = 2950.5
$ws.Range("L46").Value = 30046
$ws.Range("M46").Value = -2794.5
$ws.Range("N46").Value = -30358

$ws.Range("H102").Value = 1742.4
$ws.Range("I102").Value = 1918
$ws.Range("K102").Value = 1918
$ws.Range("M102").Value = -296

$ws.Range("H107").Value = 1907.3684
$ws.Range("I107").Value = 776.46155
$ws.Range("K107").Value = 776.46155
$ws.Range("M107").Value = 1143.53845

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H128").Value = 49250

$ws.Range("H132").Value = 3107.5881
$ws.Range("I132").Value = 2823.6428
$ws.Range("J132").Value = 4432.6665
$ws.Range("K132").Value = 8470.928400000001
$ws.Range("L132").Value = 13297.9995
$ws.Range("M132").Value = -5940.928400000001
$ws.Range("N132").Value = -18357.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 994.5
$ws.Range("I46").Value = 994.5
$ws.Range("K46").Value = 994.5
$ws.Range("M46").Value = -806.5

$ws.Range("H103").Value = 30097
$ws.Range("J103").Value = 30097
$ws.Range("L103").Value = 30097
$ws.Range("N103").Value = -32441

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 20902.6
$ws.Range("J105").Value = 20902.6
$ws.Range("L105").Value = 20902.6
$ws.Range("N105").Value = -27890.6

$ws.Range("H126").Value = 1954.7778
$ws.Range("I126").Value = 1954.7778
$ws.Range("K126").Value = 5864.3334
$ws.Range("M126").Value = -3394.3334

$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960

$ws.Range("H136").Value = 1753.8572
$ws.Range("I136").Value = 1753.8572
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5261.571599999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2711.571599999999
$ws.Range("N136").ClearContents()
